$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update test data strings - rename "18sep" dated values to "30sep"
$ws.Range("A5").Value = "TestAuto_POC30sep"
$ws.Range("B5").Value = "TestAuto_POC30sep"
$ws.Range("C5").Value = "Facility_POC30sep"
$ws.Range("D5").Value = "Facility_POC30sep"
$ws.Range("E5").Value = "Pharmacy_POC30sep"
$ws.Range("F5").Value = "Pharmacy_POC30sep"
$ws.Range("H5").Value = "AlignmentProjectPOC30sep"
